$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vegfc"
$ws.Cells.Item(2, 3).Value = "Nrp2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.164278333333333
$ws.Cells.Item(2, 8).Value = 3.492835
$ws.Cells.Item(2, 9).Value = 0.2070506538112546
$ws.Cells.Item(2, 10).Value = 0.2070506538112546
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 26.83081766666667
$ws.Cells.Item(2, 14).Value = 80.49245300000001
$ws.Cells.Item(2, 15).Value = 0.5916656861001716
$ws.Cells.Item(2, 16).Value = 0.5916656861001716
$ws.Cells.Item(2, 17).Value = 31.23853967491723
$ws.Cells.Item(2, 18).Value = 281.146857074255
$ws.Cells.Item(2, 19).Value = 0.1225047671447251
$ws.Cells.Item(2, 20).Value = 0.1225047671447251

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vegfc"
$ws.Cells.Item(3, 3).Value = "Nrp2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.164278333333333
$ws.Cells.Item(3, 8).Value = 3.492835
$ws.Cells.Item(3, 9).Value = 0.2070506538112546
$ws.Cells.Item(3, 10).Value = 0.2070506538112546
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 10.21969166666667
$ws.Cells.Item(3, 14).Value = 30.659075
$ws.Cells.Item(3, 15).Value = 0.2253617819930474
$ws.Cells.Item(3, 16).Value = 0.2253617819930474
$ws.Cells.Item(3, 17).Value = 11.89856558084722
$ws.Cells.Item(3, 18).Value = 107.087090227625
$ws.Cells.Item(3, 19).Value = 0.0466613043057299
$ws.Cells.Item(3, 20).Value = 0.0466613043057299

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vegfc"
$ws.Cells.Item(4, 3).Value = "Nrp2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.164278333333333
$ws.Cells.Item(4, 8).Value = 3.492835
$ws.Cells.Item(4, 9).Value = 0.2070506538112546
$ws.Cells.Item(4, 10).Value = 0.2070506538112546
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 8.297426666666667
$ws.Cells.Item(4, 14).Value = 24.89228
$ws.Cells.Item(4, 15).Value = 0.1829725319067811
$ws.Cells.Item(4, 16).Value = 0.1829725319067811
$ws.Cells.Item(4, 17).Value = 9.660514090422222
$ws.Cells.Item(4, 18).Value = 86.9446268138
$ws.Cells.Item(4, 19).Value = 0.03788458236079967
$ws.Cells.Item(4, 20).Value = 0.03788458236079967

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vegfc"
$ws.Cells.Item(5, 3).Value = "Nrp2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.886460333333333
$ws.Cells.Item(5, 8).Value = 11.659381
$ws.Cells.Item(5, 9).Value = 0.6911527338349851
$ws.Cells.Item(5, 10).Value = 0.6911527338349851
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 26.83081766666667
$ws.Cells.Item(5, 14).Value = 80.49245300000001
$ws.Cells.Item(5, 15).Value = 0.5916656861001716
$ws.Cells.Item(5, 16).Value = 0.5916656861001716
$ws.Cells.Item(5, 17).Value = 104.2769085723993
$ws.Cells.Item(5, 18).Value = 938.4921771515931
$ws.Cells.Item(5, 19).Value = 0.4089313564644857
$ws.Cells.Item(5, 20).Value = 0.4089313564644857

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vegfc"
$ws.Cells.Item(6, 3).Value = "Nrp2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 3.886460333333333
$ws.Cells.Item(6, 8).Value = 11.659381
$ws.Cells.Item(6, 9).Value = 0.6911527338349851
$ws.Cells.Item(6, 10).Value = 0.6911527338349851
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 10.21969166666667
$ws.Cells.Item(6, 14).Value = 30.659075
$ws.Cells.Item(6, 15).Value = 0.2253617819930474
$ws.Cells.Item(6, 16).Value = 0.2253617819930474
$ws.Cells.Item(6, 17).Value = 39.71842628139723
$ws.Cells.Item(6, 18).Value = 357.465836532575
$ws.Cells.Item(6, 19).Value = 0.1557594117264186
$ws.Cells.Item(6, 20).Value = 0.1557594117264186

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vegfc"
$ws.Cells.Item(7, 3).Value = "Nrp2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 3.886460333333333
$ws.Cells.Item(7, 8).Value = 11.659381
$ws.Cells.Item(7, 9).Value = 0.6911527338349851
$ws.Cells.Item(7, 10).Value = 0.6911527338349851
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 8.297426666666667
$ws.Cells.Item(7, 14).Value = 24.89228
$ws.Cells.Item(7, 15).Value = 0.1829725319067811
$ws.Cells.Item(7, 16).Value = 0.1829725319067811
$ws.Cells.Item(7, 17).Value = 32.24761960874222
$ws.Cells.Item(7, 18).Value = 290.22857647868
$ws.Cells.Item(7, 19).Value = 0.1264619656440808
$ws.Cells.Item(7, 20).Value = 0.1264619656440808

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Vegfc"
$ws.Cells.Item(8, 3).Value = "Nrp2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5724183333333334
$ws.Cells.Item(8, 8).Value = 1.717255
$ws.Cells.Item(8, 9).Value = 0.1017966123537602
$ws.Cells.Item(8, 10).Value = 0.1017966123537602
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 26.83081766666667
$ws.Cells.Item(8, 14).Value = 80.49245300000001
$ws.Cells.Item(8, 15).Value = 0.5916656861001716
$ws.Cells.Item(8, 16).Value = 0.5916656861001716
$ws.Cells.Item(8, 17).Value = 15.35845193072389
$ws.Cells.Item(8, 18).Value = 138.226067376515
$ws.Cells.Item(8, 19).Value = 0.06022956249096075
$ws.Cells.Item(8, 20).Value = 0.06022956249096074

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Vegfc"
$ws.Cells.Item(9, 3).Value = "Nrp2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5724183333333334
$ws.Cells.Item(9, 8).Value = 1.717255
$ws.Cells.Item(9, 9).Value = 0.1017966123537602
$ws.Cells.Item(9, 10).Value = 0.1017966123537602
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 10.21969166666667
$ws.Cells.Item(9, 14).Value = 30.659075
$ws.Cells.Item(9, 15).Value = 0.2253617819930474
$ws.Cells.Item(9, 16).Value = 0.2253617819930474
$ws.Cells.Item(9, 17).Value = 5.84993887101389
$ws.Cells.Item(9, 18).Value = 52.649449839125
$ws.Cells.Item(9, 19).Value = 0.02294106596089887
$ws.Cells.Item(9, 20).Value = 0.02294106596089887

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Vegfc"
$ws.Cells.Item(10, 3).Value = "Nrp2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5724183333333334
$ws.Cells.Item(10, 8).Value = 1.717255
$ws.Cells.Item(10, 9).Value = 0.1017966123537602
$ws.Cells.Item(10, 10).Value = 0.1017966123537602
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 8.297426666666667
$ws.Cells.Item(10, 14).Value = 24.89228
$ws.Cells.Item(10, 15).Value = 0.1829725319067811
$ws.Cells.Item(10, 16).Value = 0.1829725319067811
$ws.Cells.Item(10, 17).Value = 4.749599143488889
$ws.Cells.Item(10, 18).Value = 42.7463922914
$ws.Cells.Item(10, 19).Value = 0.01862598390190062
$ws.Cells.Item(10, 20).Value = 0.01862598390190062
